$d = $word.ActiveDocument

# --- Locate the paragraphs we need to touch by their text content, so the
#     script does not depend on brittle positional assumptions. ---
$titleText = "Play Daltanious Free: A Nostalgic Anime Slot Machine"
$oldBlurbText = "Experience the excitement of 80s Japanese anime with the Daltanious online slot machine. Play free, win big with unique symbols, and numerous bonus features!"

$titleParaIndex = $null
$dupTitleParaIndex = $null
$blurbParaIndex = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $titleText) {
        if ($null -eq $titleParaIndex) {
            $titleParaIndex = $i
        } else {
            $dupTitleParaIndex = $i
        }
    } elseif ($text -eq $oldBlurbText) {
        $blurbParaIndex = $i
    }
}

# --- 1. Insert a new "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs.Item($titleParaIndex)
$titlePara.Range.InsertParagraphAfter() | Out-Null
$metaPara = $d.Paragraphs.Item($titleParaIndex + 1)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience the excitement of 80s Japanese anime with the Daltanious online slot machine. Play free, win big with unique symbols, and numerous bonus features!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# The new paragraph shifted every later paragraph index up by one.
$dupTitleParaIndex = $dupTitleParaIndex + 1
$blurbParaIndex = $blurbParaIndex + 1

# --- 2. Remove the duplicated "Play Daltanious Free..." paragraph near the end ---
$dupTitlePara = $d.Paragraphs.Item($dupTitleParaIndex)
$dupTitlePara.Range.Delete()

# Deleting that paragraph shifted every later paragraph index down by one.
$blurbParaIndex = $blurbParaIndex - 1

# --- 3. Replace the old "Experience the excitement..." blurb with the new prompt text ---
$blurbPara = $d.Paragraphs.Item($blurbParaIndex)
$blurbPara.Range.Find.Execute($oldBlurbText, $true, $false, $false, $false, $false, $true, 1, $false, "Prompt: Create a feature image for Daltanious that features a happy Maya warrior with glasses in a cartoon-style. The warrior should be holding a slot machine lever, and there should be an anime robot in the background. The overall vibe of the image should be fun and nostalgic, reminiscent of 80s anime.", 2) | Out-Null
